# Work against the live ActiveWorkbook/ActiveSheet exposed by the host.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets("103_2")

# The "Start Date" row (row 4, label in A4) held the wrong start date for
# Senate/House/Total (B4:D4) -- serial 34357 = 1/23/1994. Correct it to
# 1/25/1994 (serial 34359) across all three columns.
$ws.Range("B4:D4").Value = 34359

$wb.Save()
